$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Remove the ADDITIONAL_DESCRIPTION (column E) values for the ARV_REGIMEN_ADULT rows (63-131);
# the data dictionary no longer breaks these out by 1st/2nd/3rd line grouping text.
$ws.Range("E63:E131").ClearContents()

# Append the new ARV_REGIMEN_CHILD concept rows (132-158) that enumerate all regimens on the NDR.
$ws.Cells.Item(132, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(132, 3).Value = '1b'
$ws.Cells.Item(132, 4).Value = 'AZT-3TC-NVP'
$ws.Cells.Item(133, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(133, 3).Value = '1z'
$ws.Cells.Item(133, 4).Value = 'AZT-3TC-DTG'
$ws.Cells.Item(134, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(134, 3).Value = '1a'
$ws.Cells.Item(134, 4).Value = 'AZT-3TC-EFV'
$ws.Cells.Item(135, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(135, 3).Value = '1q'
$ws.Cells.Item(135, 4).Value = 'D4T-3TC-EFV'
$ws.Cells.Item(136, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(136, 3).Value = '1m'
$ws.Cells.Item(136, 4).Value = 'TDF-3TC-DTG'
$ws.Cells.Item(137, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(137, 3).Value = '1x'
$ws.Cells.Item(137, 4).Value = 'DDI-3TC-EFV'
$ws.Cells.Item(138, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(138, 3).Value = '1e'
$ws.Cells.Item(138, 4).Value = 'TDF-3TC-EFV'
$ws.Cells.Item(139, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(139, 3).Value = '1c'
$ws.Cells.Item(139, 4).Value = 'TDF-FTC-EFV'
$ws.Cells.Item(140, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(140, 3).Value = '1o'
$ws.Cells.Item(140, 4).Value = 'ABC-3TC-DTG'
$ws.Cells.Item(141, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(141, 3).Value = '4c'
$ws.Cells.Item(141, 4).Value = 'ABC-3TC-EFV'
$ws.Cells.Item(142, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(142, 3).Value = '4d'
$ws.Cells.Item(142, 4).Value = 'ABC-3TC-NVP'
$ws.Cells.Item(143, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(143, 3).Value = '4f'
$ws.Cells.Item(143, 4).Value = 'D4T-3TC-NVP'
$ws.Cells.Item(144, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(144, 3).Value = '4g'
$ws.Cells.Item(144, 4).Value = 'Child First Line Others'
$ws.Cells.Item(145, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(145, 3).Value = '5a'
$ws.Cells.Item(145, 4).Value = 'ABC-3TC-LPV/r'
$ws.Cells.Item(146, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(146, 3).Value = '2b'
$ws.Cells.Item(146, 4).Value = 'TDF-3TC-LPV/r'
$ws.Cells.Item(147, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(147, 3).Value = '2a'
$ws.Cells.Item(147, 4).Value = 'TDF-FTC-LPV/r'
$ws.Cells.Item(148, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(148, 3).Value = '2e'
$ws.Cells.Item(148, 4).Value = 'AZT-3TC-LPV/r'
$ws.Cells.Item(149, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(149, 3).Value = '5c'
$ws.Cells.Item(149, 4).Value = 'D4T-3TC-LPV/r'
$ws.Cells.Item(150, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(150, 3).Value = '5d'
$ws.Cells.Item(150, 4).Value = 'DDI-3TC-NVP'
$ws.Cells.Item(151, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(151, 3).Value = '5f'
$ws.Cells.Item(151, 4).Value = 'Child Second Line Others'
$ws.Cells.Item(152, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(152, 3).Value = '6a'
$ws.Cells.Item(152, 4).Value = 'DRV/r + 2 NRTIs + 2 NNRTI'
$ws.Cells.Item(153, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(153, 3).Value = '6b'
$ws.Cells.Item(153, 4).Value = 'DRV/r +2NRTIs'
$ws.Cells.Item(154, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(154, 3).Value = '6c'
$ws.Cells.Item(154, 4).Value = 'DRV/r-DTG + 1-2 NRTIs'
$ws.Cells.Item(155, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(155, 3).Value = '6d'
$ws.Cells.Item(155, 4).Value = 'DRV/r+RAL + 1-2NRTIs'
$ws.Cells.Item(156, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(156, 3).Value = '6e'
$ws.Cells.Item(156, 4).Value = 'DTG+2 NRTIs'
$ws.Cells.Item(157, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(157, 3).Value = '6f'
$ws.Cells.Item(157, 4).Value = 'RAL(or DTG) + 2 NRTIs'
$ws.Cells.Item(158, 2).Value = 'ARV_REGIMEN_CHILD'
$ws.Cells.Item(158, 3).Value = '6g'
$ws.Cells.Item(158, 4).Value = 'DRV/r-2NRTIs+NNRTI'

# Restore the view/selection state that Excel records after making this edit.
$ws.Activate()
$ws.Range("B63:D158").Select()
